# Added Week 15 simulations
# Updates row 2 (H) target-depth counts on both the OFF and DEF sheets
# to reflect the additional week of simulated data.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 403
$wsOff.Range("C2").Value = 300
$wsOff.Range("D2").Value = 103
$wsOff.Range("E2").Value = 48
$wsOff.Range("F2").Value = 10

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 409
$wsDef.Range("C2").Value = 286
$wsDef.Range("D2").Value = 80
$wsDef.Range("E2").Value = 31
